$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.862
$ws.Range("B4").Value = 7.014999999999999
$ws.Range("D4").Value = -7.542999999999999

$ws.Range("B5").Value = 6.226

$ws.Range("A7").Value = -20.987

$ws.Range("B8").Value = 6.275999999999999

$ws.Range("D9").Value = -7.868

$ws.Range("A16").Value = -20.931
$ws.Range("B16").Value = 6.436

$ws.Range("D18").Value = -8.43
